# Insert a new data row at row 755 (pushing the existing rows 755-814 down
# to 756-815) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(755).Insert()

$ws.Range("A755").Value = 10
$ws.Range("B755").Value = "Vega Modelo de Temuco"
$ws.Range("C755").Value = "La Araucanía"
$ws.Range("D755").Value = 45106
$ws.Range("E755").Value = 9
$ws.Range("F755").Value = 100112043
$ws.Range("G755").Value = "Pepino ensalada"
$ws.Range("H755").Value = "Sin especificar"
$ws.Range("I755").Value = "Primera"
$ws.Range("J755").Value = 255
$ws.Range("K755").Value = 17000
$ws.Range("L755").Value = 17000
$ws.Range("M755").Value = 17000
$ws.Range("N755").Value = "$/caja 60 unidades"
$ws.Range("O755").Value = "Región de Arica y Parinacota"
$ws.Range("P755").Value = 283
$ws.Range("Q755").Value = 60
$ws.Range("R755").Value = "Hortaliza"
